# Commit: "Add final report + pres"
#
# Slide 16 ("Future Work (1)") body placeholder had a bullet point
# "Incomplete support in current Tomcat-Native" under "Implement TLS
# Sessions" that is no longer needed — remove that whole paragraph so
# "Was not important for a prototype" immediately follows "Implement TLS
# Sessions".

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(16)
$shp = $s.Shapes.Item(3)          # body placeholder (idx=1, type=body)
$tr  = $shp.TextFrame.TextRange

# Paragraphs are 1-based; paragraph 4 is
# "Incomplete support in current Tomcat-Native".
$target = $tr.Paragraphs(4, 1)
$target.Delete()
